$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Наименование", "Код товара", "Цена", "Ссылка на картинку"),
    @("5`" Смартфон bright & quick BQ 5016G CHOICE 16 ГБ красный", "1646107", "3999", "https://s.technopoint.ru/thumb/st1/fit/800/650/853093d89e3ff3610be956c1fd95c3e2/5c126190673c5232f844241ac215a576c667e82acb9852679508270d9ae99038.jpg"),
    @("5`" Смартфон bright & quick BQ 5016G CHOICE 16 ГБ синий", "1646102", "3999", "None"),
    @("5`" Смартфон bright & quick BQ 5016G CHOICE 16 ГБ черный", "1646104", "3999", "None"),
    @("5.99`" Смартфон TP-Link Neffos A5 16 ГБ зеленый", "1645574", "3999", "https://s.technopoint.ru/thumb/st1/fit/800/650/9710105c98ee88b56ff0a7d79e1e8a05/6557ae69ab8ab3b7c6fc2918952a6d7105c54c7c0ea40dfc6b2dcaef45b1dbdb.png"),
    @("5.99`" Смартфон TP-Link Neffos A5 16 ГБ серый", "1645572", "3999", "https://s.technopoint.ru/thumb/st1/fit/800/650/06e76fe70d7e293fb5e8627c2a4bdfe8/c958b42f59f71fce27c8bbd954be0361804ce47e23933231a127f45e98cc9a9e.png"),
    @("5.99`" Смартфон TP-Link Neffos A5 16 ГБ фиолетовый", "1645575", "3999", "https://s.technopoint.ru/thumb/st4/fit/800/650/cd8dbea757059d26fc62dbc557f1b0d5/86d2c79d23d1def8eaa4ff6b126897a23d823c3c75448860d698819714df281e.png"),
    @("6.09`" Смартфон bright & quick BQ 6042L MAGIC E 32 ГБ красный", "1646114", "5399", "https://s.technopoint.ru/thumb/st4/fit/800/650/44eb27efde71257fd8cad2194c5ef83f/d82c3ec8a5cd2baccd2a464cea4a385abfc879cfbba7fb3d37d57564182fc699.jpg"),
    @("6.09`" Смартфон bright & quick BQ 6042L MAGIC E 32 ГБ синий", "1646112", "5399", "https://s.technopoint.ru/thumb/st4/fit/800/650/44eb27efde71257fd8cad2194c5ef83f/d82c3ec8a5cd2baccd2a464cea4a385abfc879cfbba7fb3d37d57564182fc699.jpg"),
    @("6.09`" Смартфон bright & quick BQ 6042L MAGIC E 32 ГБ фиолетовый", "1646113", "5999", "https://s.technopoint.ru/thumb/st4/fit/800/650/44eb27efde71257fd8cad2194c5ef83f/d82c3ec8a5cd2baccd2a464cea4a385abfc879cfbba7fb3d37d57564182fc699.jpg"),
    @("6.09`" Смартфон bright & quick BQ 6042L MAGIC E 32 ГБ черный", "1646108", "5999", "https://s.technopoint.ru/thumb/st4/fit/800/650/44eb27efde71257fd8cad2194c5ef83f/d82c3ec8a5cd2baccd2a464cea4a385abfc879cfbba7fb3d37d57564182fc699.jpg")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $row[0]
    $ws.Cells.Item($rowNum, 2).Value = $row[1]
    $ws.Cells.Item($rowNum, 3).Value = $row[2]
    $ws.Cells.Item($rowNum, 4).Value = $row[3]
}
